$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A13").Value = 9556.11
$ws.Range("B13").Value = 9543.7000000000007
$ws.Range("C13").Value = 77.78
$ws.Range("D13").Value = 77.88
$ws.Range("E13").Value = $false
$ws.Range("F13").Value = 0.13
$ws.Range("G12").Copy($ws.Range("G13"))
$ws.Range("G13").Value = 42620.766111111108
$ws.Range("H13").Value = $true
